$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill of Materials-BOM")

# Insert a new row at row 8, pushing existing rows 8-24 down to 9-25.
$ws.Rows.Item(8).Insert()

# Populate the new row with the diode data (D1 / SMAJ5.0 TVS diode).
$ws.Cells.Item(8, 1).Value = "D1"
$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 3).Value = "Taiwan Semiconductor"
$ws.Cells.Item(8, 4).Value = "SMAJ5.0"
$ws.Cells.Item(8, 6).Value = "DO-214AC-2 Taiwan Semiconductor"
$ws.Cells.Item(8, 7).Value = "400W, 6.9V, 10%, Unidirectional, TVS"
$ws.Cells.Item(8, 12).Value = "SMAJ5.0"
